# daily auto push: 2026-02-16 10:06 UTC
# Insert a new data row for 2026/02/16 17:00 just before the existing
# row 805 (2026/12/29), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 805; everything from the old row 805
# downward (previously ending at row 846) moves down to row 806..847.
$ws.Rows.Item(805).Insert()

# Fill in the new row. Column A holds a date formatted as plain text
# (e.g. "2026/02/16"), so force text formatting first to avoid Excel's
# automatic date auto-detection, then restore the default "Normal"
# style so the cell matches the unstyled look of the surrounding data
# cells once the value has been stored as text.
$dateCell = $ws.Cells.Item(805, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026/02/16"
$dateCell.Style = "Normal"

$ws.Cells.Item(805, 2).Value = "月"
$ws.Cells.Item(805, 3).Value = 17
$ws.Cells.Item(805, 4).Value = 44
